# Fruta / hortaliza, semanal
# A new weekly observation is inserted at the top of the Cilantro price
# series (row 65). All subsequent rows shift down by one, and the sheet
# dimension grows from A1:R99 to A1:R100.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 65, pushing rows 65-99 down to 66-100.
$ws.Rows.Item(65).Insert()

# Populate the newly inserted row 65 with the new weekly data point.
$ws.Range("A65").Value = 1
$ws.Range("B65").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C65").Value = "Arica y Parinacota"
$ws.Range("D65").Value = 44917
$ws.Range("E65").Value = 15
$ws.Range("F65").Value = 100112040
$ws.Range("G65").Value = "Cilantro"
$ws.Range("H65").Value = "Sin especificar"
$ws.Range("I65").Value = "Primera"
$ws.Range("J65").Value = 400
$ws.Range("K65").Value = 2700
$ws.Range("L65").Value = 3000
$ws.Range("M65").Value = 2888
$ws.Range("N65").Value = "$/atado 1,5 a 2 kilos"
$ws.Range("O65").Value = "Región de Arica y Parinacota"
$ws.Range("P65").Value = 1444
$ws.Range("Q65").Value = 2
$ws.Range("R65").Value = "Hortaliza"
